# Updated cryptos list on Tue Jun 25 15:25:40 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.722.50"
$ws.Range("E2").Value = "  +0.70%  "

$ws.Range("D3").Value = "3.416.17"
$ws.Range("E3").Value = "  +3.24%  "

$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.26%  "

$ws.Range("D5").Value = "576.86"
$ws.Range("E5").Value = "  +1.72%  "

$ws.Range("D6").Value = "138.47"

$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.21%  "

$ws.Range("D8").Value = "3.415.40"
$ws.Range("E8").Value = "  +3.28%  "

$ws.Range("D9").Value = "0.476"
$ws.Range("E9").Value = "  -0.15%  "

$ws.Range("D10").Value = "7.68"
$ws.Range("E10").Value = "  +4.58%  "

$ws.Range("E11").Value = "  +5.91%  "

$ws.Range("D12").Value = "0.395"
$ws.Range("E12").Value = "  +5.31%  "

$ws.Range("D13").Value = "3.980.64"
$ws.Range("E13").Value = "  +2.85%  "

$ws.Range("E14").Value = "  +2.13%  "

$ws.Range("D15").Value = "0.0000179"
$ws.Range("E15").Value = "  +6.83%  "

$ws.Range("D16").Value = "3.414.27"
$ws.Range("E16").Value = "  +3.38%  "

$ws.Range("D17").Value = "25.54"
$ws.Range("E17").Value = "  +5.03%  "

$ws.Range("D18").Value = "61.602.64"
$ws.Range("E18").Value = "  +0.51%  "

$ws.Range("D19").Value = "14.08"
$ws.Range("E19").Value = "  +4.97%  "

$ws.Range("D20").Value = "5.92"
$ws.Range("E20").Value = "  +4.45%  "

$ws.Range("D21").Value = "9.47"
$ws.Range("E21").Value = "  +5.68%  "

$ws.Range("D22").Value = "389.00"
$ws.Range("E22").Value = "  +9.66%  "

$ws.Range("D23").Value = "0.574"
$ws.Range("E23").Value = "  +3.57%  "

$ws.Range("D24").Value = "3.536.80"
$ws.Range("E24").Value = "  +2.84%  "

$ws.Range("E25").Value = "  +0.11%  "

$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").Value = "0.0000127"
$ws.Range("E26").Value = "  +18.23%  "

$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").Value = "71.02"
$ws.Range("E27").Value = "  +2.82%  "

$ws.Range("D28").Value = "1.67"
$ws.Range("E28").Value = "  +15.78%  "

$ws.Range("D29").Value = "7.83"
$ws.Range("E29").Value = "  +9.16%  "

$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.46%  "

$ws.Range("D31").Value = "8.33"
$ws.Range("E31").Value = "  +6.61%  "

$ws.Range("E32").Value = "  +6.02%  "

$ws.Range("E33").Value = "  +2.35%  "

$ws.Range("B34").Value = "RenzoRestakedETH"
$ws.Range("C34").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D34").Value = "3.438.22"
$ws.Range("E34").Value = "  +3.03%  "

$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("D36").Value = "23.65"
$ws.Range("E36").Value = "  +4.53%  "

$ws.Range("D37").Value = "5.54"
$ws.Range("E37").Value = "  +4.49%  "

$ws.Range("E38").Value = "  +3.89%  "

$ws.Range("E39").Value = "  +6.00%  "

$ws.Range("D40").Value = "162.11"
$ws.Range("E40").Value = "  -0.44%  "

$ws.Range("D41").Value = "0.0799"
$ws.Range("E41").Value = "  +5.77%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "1.73"
$ws.Range("E42").Value = "  +11.42%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "0.996"
$ws.Range("E43").Value = "  -0.25%  "

$ws.Range("D44").Value = "1.23"
$ws.Range("E44").Value = "  +9.01%  "

$ws.Range("D45").Value = "0.774"
$ws.Range("E45").Value = "  +4.18%  "

$ws.Range("D46").Value = "4.46"
$ws.Range("E46").Value = "  +1.42%  "

$ws.Range("D47").Value = "41.31"
$ws.Range("E47").Value = "  +0.48%  "

$ws.Range("D48").Value = "23.52"
$ws.Range("E48").Value = "  +5.08%  "

$ws.Range("D49").Value = "7.03"
$ws.Range("E49").Value = "  +5.27%  "

$ws.Range("D50").Value = "22.99"
$ws.Range("E50").Value = "  +7.93%  "

$ws.Range("D51").Value = "2.352.17"
$ws.Range("E51").Value = "  +8.80%  "
